$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: task names ---
$ws.Range("A1").Value = "Btn handler task"
$ws.Range("A2").Value = "Platform diraction task"
$ws.Range("A3").Value = "physics task"
$ws.Range("A4").Value = "LCD display task"
$ws.Range("A5").Value = "LED display task"

# --- Column B: durations (all "1 week" now) ---
$ws.Range("B1").Value = "1 week"
$ws.Range("B2").Value = "1 week"
$ws.Range("B3").Value = "1 week"
$ws.Range("B4").Value = "1 week"
$ws.Range("B5").Value = "1 week"

# --- Column C: new status column ---
$ws.Range("C1").Value = "done"
$ws.Range("C2").Value = "not done"
$ws.Range("C3").Value = "not done"
$ws.Range("C4").Value = "not done"
$ws.Range("C5").Value = "not done"

# Match the existing table formatting (style index used by columns A/B)
# on the freshly-written column C cells.
$ws.Range("B1:B5").Copy()
$ws.Range("C1:C5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
